# Apply the row permutation described by the diff to rows 12-17.
# The content (all columns A:AY) of these six rows gets permuted as follows
# (after[row] = before[mapping[row]]):
#   12 <- 17
#   13 <- 16
#   14 <- 13
#   15 <- 12
#   16 <- 15
#   17 <- 14
# i.e. a single 6-cycle: 12 -> 17 -> 14 -> 13 -> 16 -> 15 -> 12
#
# All the cell values in these rows are plain text/number/boolean data (no
# formulas), so the permutation is implemented by reading each row's full
# A:AY values, then writing those values back into the new row positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of text cells in this sheet contain digit-only or ISO-date-like
# text (e.g. Antal="1", Startdatum="2023-08-28") even though the column as a
# whole is textual. When such a string is written back through Value2,
# Excel auto-converts it into a real number/date, which would change both
# the stored value and its cell type. Prefixing with a leading apostrophe
# forces Excel to keep (and store) the value as plain text, matching the
# original representation, without leaving the apostrophe in the value.
function Protect-TextValue {
    param($v)
    if ($v -is [string] -and $v.Length -gt 0) {
        if ($v -match '^\s*-?\d+(\.\d+)?\s*$' -or $v -match '^\s*\d{4}-\d{2}-\d{2}\s*$') {
            return "'" + $v
        }
    }
    return $v
}

function Get-ProtectedRow {
    param($range)
    $vals = $range.Value2
    for ($c = 1; $c -le $vals.GetLength(1); $c++) {
        $vals[1, $c] = Protect-TextValue $vals[1, $c]
    }
    return $vals
}

# Capture the original contents of each row first, before any writes, so
# that source data is not clobbered mid-way through the permutation.
$row12 = Get-ProtectedRow $ws.Range("A12:AY12")
$row13 = Get-ProtectedRow $ws.Range("A13:AY13")
$row14 = Get-ProtectedRow $ws.Range("A14:AY14")
$row15 = Get-ProtectedRow $ws.Range("A15:AY15")
$row16 = Get-ProtectedRow $ws.Range("A16:AY16")
$row17 = Get-ProtectedRow $ws.Range("A17:AY17")

# Write back according to the permutation.
$ws.Range("A12:AY12").Value2 = $row17
$ws.Range("A13:AY13").Value2 = $row16
$ws.Range("A14:AY14").Value2 = $row13
$ws.Range("A15:AY15").Value2 = $row12
$ws.Range("A16:AY16").Value2 = $row15
$ws.Range("A17:AY17").Value2 = $row14
